# Apply the backlog update: rename "arrousi" -> "Omar Al Arousi" and add a
# new "I can delete events" user story row (row 53) assigned to her, on
# Sheet1 of the backlog workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the new user story row (row 53) ------------------------------
# Set B53 first so the new "I can delete events " string is appended to the
# shared string table before the "Omar Al Arousi" rename below, matching
# the order strings were added upstream.
$ws.Range("A53").Value = 2.13
$ws.Range("B53").Value = "I can delete events "
$ws.Range("C53").Value = 8

# --- Rename the assignee "arrousi" -> "Omar Al Arousi" --------------------
# This touches every cell in column D that previously held "arrousi"
# (rows 48-52), plus the new row 53 below.
$ws.Range("D48").Value = "Omar Al Arousi"
$ws.Range("D49").Value = "Omar Al Arousi"
$ws.Range("D50").Value = "Omar Al Arousi"
$ws.Range("D51").Value = "Omar Al Arousi"
$ws.Range("D52").Value = "Omar Al Arousi"
$ws.Range("D53").Value = "Omar Al Arousi"

# E53 uses the same "0.00" number style as the other "Depends on" values in
# this block (e.g. E51/E52), not the row's default style.
$ws.Range("E53").Value = 2.1
$ws.Range("E53").NumberFormat = "0.00"

# --- Update the view/selection state --------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 37
$ws.Range("E55").Select()
